$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape "Retangulo 68" (id=69) -> Shapes.Item(21)
# off: 6493817,3775025 -> 6493816,3775025 ; ext: 1375356,600163 -> 1385769,768380 (EMU)
$sh1 = $s.Shapes.Item(21)
$sh1.Left = 511.32409448818896
$sh1.Top = 297.246062992126
$sh1.Width = 109.11567309133858
$sh1.Height = 60.502363304724405

# Shape "CaixaDeTexto 70" (id=71) -> Shapes.Item(22)
# off: 6451330,3800440 -> 6451330,3800441 (x unchanged) ; ext: 1411806,600164 -> 1385781,769441 (EMU)
$sh2 = $s.Shapes.Item(22)
$sh2.Top = 299.2473297346457
$sh2.Width = 109.11661417322834
$sh2.Height = 60.58590551181102
$sh2.TextFrame.TextRange.Text = "Subsistema de confirmação do pedido (carrinho de compra)"

# Shape "CaixaDeTexto 193" (id=194) -> Shapes.Item(54)
# off: 4427022,1743016 -> 4368354,1743016 (y unchanged)
$sh3 = $s.Shapes.Item(54)
$sh3.Left = 343.9648895897638
$sh3.TextFrame.TextRange.Text = "Manter cadastros"
